$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value (serial 45205 = 2023-10-06) for every
# data row (rows 2 through 303). Update it to 45206 (2023-10-07) for all of them.
for ($row = 2; $row -le 303; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
